$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking id columns stay stored as text (matching source inlineStr type)
$ws.Range("G1:G2").NumberFormat = "@"
$ws.Range("P1:P2").NumberFormat = "@"

# Row 1 changes
$ws.Range("C1").Value = '[{''hashtag_description'': '''', ''hashtag_id'': 45455, ''hashtag_name'': ''chickenwings''}, {''hashtag_description'': '''', ''hashtag_id'': 3979591, ''hashtag_name'': ''firewings''}, {''hashtag_description'': '''', ''hashtag_id'': 47867, ''hashtag_name'': ''chickenwing''}, {''hashtag_description'': '''', ''hashtag_id'': 2556710, ''hashtag_name'': ''familygames''}, {''hashtag_description'': "Whether it''s desktop, table top, or old school, what are you playing for #FamilyGameNight?", ''hashtag_id'': 287421, ''hashtag_name'': ''FamilyGameNight''}, {''hashtag_description'': ''Asla bir sonraki hamleni bilmelerine izin verme! 👀'', ''hashtag_id'': 1652484531221509, ''hashtag_name'': ''xyzbca''}]'
$ws.Range("G1").Value = '252133'
$ws.Range("N1").Value = '{''vote'': False, ''warn'': False, ''content'': '''', ''sink'': False, ''type'': 0}'
$ws.Range("P1").Value = '4314291'

# Row 2 changes
$ws.Range("C2").Value = '[{''hashtag_name'': ''xyzbca'', ''hashtag_description'': ''Asla bir sonraki hamleni bilmelerine izin verme! 👀'', ''hashtag_id'': 1652484531221509}, {''hashtag_description'': '''', ''hashtag_id'': 45455, ''hashtag_name'': ''chickenwings''}, {''hashtag_name'': ''firewings'', ''hashtag_description'': '''', ''hashtag_id'': 3979591}, {''hashtag_description'': '''', ''hashtag_id'': 47867, ''hashtag_name'': ''chickenwing''}, {''hashtag_id'': 2556710, ''hashtag_name'': ''familygames'', ''hashtag_description'': ''''}, {''hashtag_description'': "Whether it''s desktop, table top, or old school, what are you playing for #FamilyGameNight?", ''hashtag_id'': 287421, ''hashtag_name'': ''FamilyGameNight''}]'
$ws.Range("G2").Value = '71028'
$ws.Range("N2").Value = '{''content'': '''', ''sink'': False, ''type'': 0, ''vote'': False, ''warn'': False}'
$ws.Range("P2").Value = '1719544'
